$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for rows 2-25 (case with 380 kV)

# Row 2
$ws.Cells.Item(2, 2).Value = 16.87792829883318
$ws.Cells.Item(2, 4).Value = 3.756858328303573
$ws.Cells.Item(2, 5).Value = 23.6275167487601
$ws.Cells.Item(2, 6).Value = 26.62715821472407
$ws.Cells.Item(2, 7).Value = 35.17175490361991
$ws.Cells.Item(2, 8).Value = 14.28575027185112
$ws.Cells.Item(2, 12).Value = 9.071561772022136
$ws.Cells.Item(2, 13).Value = 15.50722478689785
$ws.Cells.Item(2, 14).Value = 19.41956272620733

# Row 3
$ws.Cells.Item(3, 2).Value = 16.69055607013022
$ws.Cells.Item(3, 4).Value = 3.779977242251641
$ws.Cells.Item(3, 5).Value = 23.19863771392944
$ws.Cells.Item(3, 6).Value = 26.07889045076944
$ws.Cells.Item(3, 7).Value = 34.00925087887191
$ws.Cells.Item(3, 8).Value = 14.1989755036751
$ws.Cells.Item(3, 12).Value = 8.949398383956426
$ws.Cells.Item(3, 13).Value = 15.40731268961109
$ws.Cells.Item(3, 14).Value = 19.46719898768509

# Row 4
$ws.Cells.Item(4, 2).Value = 16.57794528369668
$ws.Cells.Item(4, 4).Value = 3.795028221844162
$ws.Cells.Item(4, 5).Value = 22.92940124251381
$ws.Cells.Item(4, 6).Value = 25.74722625757903
$ws.Cells.Item(4, 7).Value = 33.2904314911795
$ws.Cells.Item(4, 8).Value = 14.15012927755302
$ws.Cells.Item(4, 12).Value = 8.875157617447242
$ws.Cells.Item(4, 13).Value = 15.34834493926704
$ws.Cells.Item(4, 14).Value = 19.49886873922952

# Row 5
$ws.Cells.Item(5, 2).Value = 16.53271735357541
$ws.Cells.Item(5, 4).Value = 3.801376346578175
$ws.Cells.Item(5, 5).Value = 22.81828687927028
$ws.Cells.Item(5, 6).Value = 25.61353594488903
$ws.Cells.Item(5, 7).Value = 32.99679117106934
$ws.Cells.Item(5, 8).Value = 14.13135360721073
$ws.Cells.Item(5, 12).Value = 8.845130448851709
$ws.Cells.Item(5, 13).Value = 15.32493245312138
$ws.Cells.Item(5, 14).Value = 19.51238429861887

# Row 6
$ws.Cells.Item(6, 2).Value = 16.52524867899719
$ws.Cells.Item(6, 4).Value = 3.802443403263049
$ws.Cells.Item(6, 5).Value = 22.79975454858593
$ws.Cells.Item(6, 6).Value = 25.59143154754616
$ws.Cells.Item(6, 7).Value = 32.94800596744622
$ws.Cells.Item(6, 8).Value = 14.12830453029629
$ws.Cells.Item(6, 12).Value = 8.840159113528784
$ws.Cells.Item(6, 13).Value = 15.32108265684287
$ws.Cells.Item(6, 14).Value = 19.51466541527971

# Row 7
$ws.Cells.Item(7, 2).Value = 16.57733258009
$ws.Cells.Item(7, 4).Value = 3.795112966099341
$ws.Cells.Item(7, 5).Value = 22.92790826076658
$ws.Cells.Item(7, 6).Value = 25.74541704622534
$ws.Cells.Item(7, 7).Value = 33.28647349933218
$ws.Cells.Item(7, 8).Value = 14.14987147111078
$ws.Cells.Item(7, 12).Value = 8.874751699396031
$ws.Cells.Item(7, 13).Value = 15.34802666604488
$ws.Cells.Item(7, 14).Value = 19.49904854394493

# Row 8
$ws.Cells.Item(8, 2).Value = 16.81284653009188
$ws.Cells.Item(8, 4).Value = 3.764651792043096
$ws.Cells.Item(8, 5).Value = 23.4809149815232
$ws.Cells.Item(8, 6).Value = 26.43721180509105
$ws.Cells.Item(8, 7).Value = 34.7722887488821
$ws.Cells.Item(8, 8).Value = 14.25491951336756
$ws.Cells.Item(8, 12).Value = 9.029300550234328
$ws.Cells.Item(8, 13).Value = 15.47229320537364
$ws.Cells.Item(8, 14).Value = 19.43548592612528

# Row 9
$ws.Cells.Item(9, 2).Value = 17.29180892491921
$ws.Cells.Item(9, 4).Value = 3.711731325941509
$ws.Cells.Item(9, 5).Value = 24.5149266793265
$ws.Cells.Item(9, 6).Value = 27.82390933414564
$ws.Cells.Item(9, 7).Value = 37.62364054515269
$ws.Cells.Item(9, 8).Value = 14.49535665266157
$ws.Cells.Item(9, 12).Value = 9.33698919817423
$ws.Cells.Item(9, 13).Value = 15.73396827794742
$ws.Cells.Item(9, 14).Value = 19.32999925494513

# Row 10
$ws.Cells.Item(10, 2).Value = 17.65112984217508
$ws.Cells.Item(10, 4).Value = 3.677037507325806
$ws.Cells.Item(10, 5).Value = 25.23933339631282
$ws.Cells.Item(10, 6).Value = 28.84875868967573
$ws.Cells.Item(10, 7).Value = 39.65338110208364
$ws.Cells.Item(10, 8).Value = 14.69187325811708
$ws.Cells.Item(10, 12).Value = 9.563895340882889
$ws.Cells.Item(10, 13).Value = 15.93594964356214
$ws.Cells.Item(10, 14).Value = 19.26411220647035

# Row 11
$ws.Cells.Item(11, 2).Value = 17.81556327193992
$ws.Cells.Item(11, 4).Value = 3.662171455564282
$ws.Cells.Item(11, 5).Value = 25.56034437617268
$ws.Cells.Item(11, 6).Value = 29.31382239886014
$ws.Cells.Item(11, 7).Value = 40.55769512294751
$ws.Cells.Item(11, 8).Value = 14.78531131535678
$ws.Cells.Item(11, 12).Value = 9.666894904159076
$ws.Cells.Item(11, 13).Value = 16.02967165338222
$ws.Cells.Item(11, 14).Value = 19.23664653332811

# Row 12
$ws.Cells.Item(12, 2).Value = 17.87791684654782
$ws.Cells.Item(12, 4).Value = 3.656674597043638
$ws.Cells.Item(12, 5).Value = 25.68060889647271
$ws.Cells.Item(12, 6).Value = 29.48957064102431
$ws.Cells.Item(12, 7).Value = 40.89704438480442
$ws.Cells.Item(12, 8).Value = 14.82124956225553
$ws.Cells.Item(12, 12).Value = 9.705833013975777
$ws.Cells.Item(12, 13).Value = 16.0654021936005
$ws.Cells.Item(12, 14).Value = 19.22660536867377

# Row 13
$ws.Cells.Item(13, 2).Value = 17.86448490977824
$ws.Cells.Item(13, 4).Value = 3.657852533131116
$ws.Cells.Item(13, 5).Value = 25.6547664515395
$ws.Cells.Item(13, 6).Value = 29.4517391994073
$ws.Cells.Item(13, 7).Value = 40.82410239766813
$ws.Cells.Item(13, 8).Value = 14.81348534702744
$ws.Cells.Item(13, 12).Value = 9.697450440047668
$ws.Cells.Item(13, 13).Value = 16.05769669968771
$ws.Cells.Item(13, 14).Value = 19.22875194005396

# Row 14
$ws.Cells.Item(14, 2).Value = 17.82069171826549
$ws.Cells.Item(14, 4).Value = 3.661716563458819
$ws.Cells.Item(14, 5).Value = 25.5702649540397
$ws.Cells.Item(14, 6).Value = 29.32828964210585
$ws.Cells.Item(14, 7).Value = 40.58567753449527
$ws.Cells.Item(14, 8).Value = 14.78825700067937
$ws.Cells.Item(14, 12).Value = 9.670099870482881
$ws.Cells.Item(14, 13).Value = 16.03260655696694
$ws.Cells.Item(14, 14).Value = 19.23581324147433

# Row 15
$ws.Cells.Item(15, 2).Value = 17.79387675379636
$ws.Cells.Item(15, 4).Value = 3.664100686370521
$ws.Cells.Item(15, 5).Value = 25.51833466487133
$ws.Cells.Item(15, 6).Value = 29.25262051675754
$ws.Cells.Item(15, 7).Value = 40.43922239128855
$ws.Cells.Item(15, 8).Value = 14.7728754340357
$ws.Cells.Item(15, 12).Value = 9.653337328892597
$ws.Cells.Item(15, 13).Value = 16.01726864887428
$ws.Cells.Item(15, 14).Value = 19.24018528007642

# Row 16
$ws.Cells.Item(16, 2).Value = 17.64039873004045
$ws.Cells.Item(16, 4).Value = 3.678027547206194
$ws.Cells.Item(16, 5).Value = 25.21817698738892
$ws.Cells.Item(16, 6).Value = 28.81832559417817
$ws.Cells.Item(16, 7).Value = 39.59386763440018
$ws.Cells.Item(16, 8).Value = 14.68584600490459
$ws.Cells.Item(16, 12).Value = 9.557156717175365
$ws.Cells.Item(16, 13).Value = 15.92985962302995
$ws.Cells.Item(16, 14).Value = 19.26595751066323

# Row 17
$ws.Cells.Item(17, 2).Value = 17.54645625538808
$ws.Cells.Item(17, 4).Value = 3.68680649286435
$ws.Cells.Item(17, 5).Value = 25.03180565209285
$ws.Cells.Item(17, 6).Value = 28.55146325385176
$ws.Cells.Item(17, 7).Value = 39.07012913754004
$ws.Cells.Item(17, 8).Value = 14.63347342037463
$ws.Cells.Item(17, 12).Value = 9.498072049375347
$ws.Cells.Item(17, 13).Value = 15.87669089078674
$ws.Cells.Item(17, 14).Value = 19.28240926835482

# Row 18
$ws.Cells.Item(18, 2).Value = 17.49251815479124
$ws.Cells.Item(18, 4).Value = 3.6919421082632
$ws.Cells.Item(18, 5).Value = 24.92381122719091
$ws.Cells.Item(18, 6).Value = 28.39787442335711
$ws.Cells.Item(18, 7).Value = 38.76711678771959
$ws.Cells.Item(18, 8).Value = 14.60373215469689
$ws.Cells.Item(18, 12).Value = 9.464069260139514
$ws.Cells.Item(18, 13).Value = 15.84628429314827
$ws.Cells.Item(18, 14).Value = 19.29210788180303

# Row 19
$ws.Cells.Item(19, 2).Value = 17.47427364217964
$ws.Cells.Item(19, 4).Value = 3.693695719661346
$ws.Cells.Item(19, 5).Value = 24.8871112043567
$ws.Cells.Item(19, 6).Value = 28.34586143895587
$ws.Cells.Item(19, 7).Value = 38.6642292837584
$ws.Cells.Item(19, 8).Value = 14.59372866030114
$ws.Cells.Item(19, 12).Value = 9.452554306814138
$ws.Cells.Item(19, 13).Value = 15.83601987846291
$ws.Cells.Item(19, 14).Value = 19.29543222804799

# Row 20
$ws.Cells.Item(20, 2).Value = 17.55644714207502
$ws.Cells.Item(20, 4).Value = 3.685863032131459
$ws.Cells.Item(20, 5).Value = 25.05172836395223
$ws.Cells.Item(20, 6).Value = 28.57988266753426
$ws.Cells.Item(20, 7).Value = 39.1260680960177
$ws.Cells.Item(20, 8).Value = 14.63900921328488
$ws.Cells.Item(20, 12).Value = 9.504363911861869
$ws.Cells.Item(20, 13).Value = 15.88233289372018
$ws.Cells.Item(20, 14).Value = 19.28063353285154

# Row 21
$ws.Cells.Item(21, 2).Value = 17.8335529289629
$ws.Cells.Item(21, 4).Value = 3.660577998747203
$ws.Cells.Item(21, 5).Value = 25.59512077999841
$ws.Cells.Item(21, 6).Value = 29.3645610896176
$ws.Cells.Item(21, 7).Value = 40.65579535800716
$ws.Cells.Item(21, 8).Value = 14.79565231994809
$ws.Cells.Item(21, 12).Value = 9.678135440288052
$ws.Cells.Item(21, 13).Value = 16.0399698220568
$ws.Cells.Item(21, 14).Value = 19.23372941730045

# Row 22
$ws.Cells.Item(22, 2).Value = 18.01513548707793
$ws.Cells.Item(22, 4).Value = 3.644826036698739
$ws.Cells.Item(22, 5).Value = 25.94267959328883
$ws.Cells.Item(22, 6).Value = 29.87520429665548
$ws.Cells.Item(22, 7).Value = 41.63738195227913
$ws.Cells.Item(22, 8).Value = 14.90125140484193
$ws.Cells.Item(22, 12).Value = 9.791307914841875
$ws.Cells.Item(22, 13).Value = 16.14438228610112
$ws.Cells.Item(22, 14).Value = 19.20516980386116

# Row 23
$ws.Cells.Item(23, 2).Value = 17.91819543209273
$ws.Cells.Item(23, 4).Value = 3.653162116043815
$ws.Cells.Item(23, 5).Value = 25.75789625627802
$ws.Cells.Item(23, 6).Value = 29.60292694675155
$ws.Cells.Item(23, 7).Value = 41.11526272955069
$ws.Cells.Item(23, 8).Value = 14.84460508157277
$ws.Cells.Item(23, 12).Value = 9.73095289006266
$ws.Cells.Item(23, 13).Value = 16.08853659384977
$ws.Cells.Item(23, 14).Value = 19.22022123765696

# Row 24
$ws.Cells.Item(24, 2).Value = 17.5519300379331
$ws.Cells.Item(24, 4).Value = 3.686289295186951
$ws.Cells.Item(24, 5).Value = 25.04272393901321
$ws.Cells.Item(24, 6).Value = 28.56703475299246
$ws.Cells.Item(24, 7).Value = 39.10078400856903
$ws.Cells.Item(24, 8).Value = 14.63650533178103
$ws.Cells.Item(24, 12).Value = 9.501519464867277
$ws.Cells.Item(24, 13).Value = 15.87978164143904
$ws.Cells.Item(24, 14).Value = 19.28143559438437

# Row 25
$ws.Cells.Item(25, 2).Value = 17.16071344047742
$ws.Cells.Item(25, 4).Value = 3.725315003905606
$ws.Cells.Item(25, 5).Value = 24.2410821882025
$ws.Cells.Item(25, 6).Value = 27.44684135777397
$ws.Cells.Item(25, 7).Value = 36.86190828961622
$ws.Cells.Item(25, 8).Value = 14.42673433259884
$ws.Cells.Item(25, 12).Value = 9.253461220268152
$ws.Cells.Item(25, 13).Value = 15.66137738678238
$ws.Cells.Item(25, 14).Value = 19.35649191341512

